$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62
$ws.Cells.Item($row, 1).Value = "Michele Bertolini 2"
$ws.Cells.Item($row, 2).Value = "Stefano Tita | Clitoriders"
$ws.Cells.Item($row, 3).Value = "Leonardo Viola | Shark Attack"
$ws.Cells.Item($row, 4).Value = "Federico  Manica | iMontagna"
$ws.Cells.Item($row, 5).Value = "Nicholas Marzadro | SBARX"
$ws.Cells.Item($row, 6).Value = "FEDERICO NICOLODI | U.S. Guarna"
